$p = $ppt.ActivePresentation

function Update-DateFld($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
        } catch {}
        if ($isDatePh -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "20/08/2024") {
                $sh.TextFrame.TextRange.Text = "06/09/2024"
            }
        }
    }
}

# Update the date placeholder on the slide master
Update-DateFld $p.SlideMaster.Shapes

# Update the date placeholder on every slide layout
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $cl = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateFld $cl.Shapes
}

# Remove the standalone banner picture from slide 1
$s = $p.Slides.Item(1)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Picture 6") {
        $sh.Delete()
    }
}
